$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet4")

# --- Fill in the new "Create Card" test case on row 3 ---
# Column order matters: new shared-string values get appended to
# sharedStrings.xml in the order they are first assigned, so we go
# left-to-right exactly like the row is laid out.
$ws.Range("A3").Value = "TS01"
$ws.Range("B3").Value = "TC002"
$ws.Range("C3").Value = "Verify that the CreateCard class creates both a Debit and Credit Card.  "
$ws.Range("D3").Value = "N/A"
$ws.Range("E3").Value = "1. User knows what kind of Card they want to create."
$ws.Range("F3").Value = "1. First Create a Instance of a Card.          2. Create a Debit Card using the instance.                                                     3. Create a Credit Card using the same instance."
$ws.Range("G3").Value = "1. New DebitCard object with name Tyler Test                 2. New CreditCard object with name Tyler Test"
$ws.Range("H3").Value = "1. First created was the Debit Card so a message that A Debit Card has been created should show.     2. The Credit Card is created after and a message that A Credit Card has been created should be displayed."
$ws.Range("I3").Value = "1. If the type of card is correctly input, the result should be as expected.                             2. If a type of card does not exist, an error message is displayed saying card cannot be created."
$ws.Range("J3").Value = "Pass"
$ws.Range("K3").Value = "Create Card Test Case"
$ws.Range("L3").Value = "Tyler Serio"

# The two "date" columns need to stay plain text (like the rest of the
# sheet) instead of being auto-recognized as serial date numbers.
$ws.Range("M3").NumberFormat = "@"
$ws.Range("M3").Value = "03/24/2015"

$ws.Range("N3").Value = "Tyler Serio"

$ws.Range("O3").NumberFormat = "@"
$ws.Range("O3").Value = "04/13/2015"

$ws.Range("P3").Value = "OS: Windows 8.1                   IDE: Eclipse"

# Row 3 matches the wrapped/tall look of row 2 (every column except the
# suite id / case id / requirement columns wraps its text).
$ws.Range("C3").WrapText = $true
$ws.Range("E3:P3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 90

# Selection moves from P4 to the newly-filled P3.
$ws.Range("P3").Select()

# Give the sheet an explicit (portrait) page setup, like the rest of
# the workbook.
$ws.PageSetup.Orientation = 1
